$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 211.4614666666667
$ws.Range("H2").Value = 634.3844
$ws.Range("I2").Value = 0.2421062275331183
$ws.Range("J2").Value = 0.2421062275331183
$ws.Range("O2").Value = 0.001498364820294181
$ws.Range("P2").Value = 0.001498364820294181
$ws.Range("Q2").Value = 2.576798945644444
$ws.Range("R2").Value = 23.1911905108
$ws.Range("S2").Value = 0.0003627634541097629
$ws.Range("T2").Value = 0.0003627634541097629
# Row 3
$ws.Range("G3").Value = 211.4614666666667
$ws.Range("H3").Value = 634.3844
$ws.Range("I3").Value = 0.2421062275331183
$ws.Range("J3").Value = 0.2421062275331183
$ws.Range("M3").Value = 0.08128566666666666
$ws.Range("N3").Value = 0.243857
$ws.Range("O3").Value = 0.009994987279658562
$ws.Range("P3").Value = 0.009994987279658561
$ws.Range("Q3").Value = 17.18878629231111
$ws.Range("R3").Value = 154.6990766308
$ws.Range("S3").Value = 0.002419848664519639
$ws.Range("T3").Value = 0.002419848664519639
# Row 4
$ws.Range("G4").Value = 211.4614666666667
$ws.Range("H4").Value = 634.3844
$ws.Range("I4").Value = 0.2421062275331183
$ws.Range("J4").Value = 0.2421062275331183
$ws.Range("M4").Value = 7.912604999999999
$ws.Range("N4").Value = 23.737815
$ws.Range("O4").Value = 0.9729438112167713
$ws.Range("P4").Value = 0.9729438112167712
$ws.Range("Q4").Value = 1673.211058454
$ws.Range("R4").Value = 15058.899526086
$ws.Range("S4").Value = 0.235555755735387
$ws.Range("T4").Value = 0.2355557557353869
# Row 5
$ws.Range("G5").Value = 211.4614666666667
$ws.Range("H5").Value = 634.3844
$ws.Range("I5").Value = 0.2421062275331183
$ws.Range("J5").Value = 0.2421062275331183
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.126567
$ws.Range("N5").Value = 0.379701
$ws.Range("O5").Value = 0.015562836683276
$ws.Range("P5").Value = 0.015562836683276
$ws.Range("Q5").Value = 26.7640434516
$ws.Range("R5").Value = 240.8763910644
$ws.Range("S5").Value = 0.00376785967910198
$ws.Range("T5").Value = 0.003767859679101979
# Row 6
$ws.Range("I6").Value = 0.08842543241393927
$ws.Range("J6").Value = 0.08842543241393927
$ws.Range("O6").Value = 0.001498364820294181
$ws.Range("P6").Value = 0.001498364820294181
$ws.Range("S6").Value = 0.0001324935571483473
$ws.Range("T6").Value = 0.0001324935571483473
# Row 7
$ws.Range("I7").Value = 0.08842543241393927
$ws.Range("J7").Value = 0.08842543241393927
$ws.Range("M7").Value = 0.08128566666666666
$ws.Range("N7").Value = 0.243857
$ws.Range("O7").Value = 0.009994987279658562
$ws.Range("P7").Value = 0.009994987279658561
$ws.Range("Q7").Value = 6.277929634670333
$ws.Range("R7").Value = 56.501366712033
$ws.Range("S7").Value = 0.0008838110721756309
$ws.Range("T7").Value = 0.0008838110721756308
# Row 8
$ws.Range("I8").Value = 0.08842543241393927
$ws.Range("J8").Value = 0.08842543241393927
$ws.Range("M8").Value = 7.912604999999999
$ws.Range("N8").Value = 23.737815
$ws.Range("O8").Value = 0.9729438112167713
$ws.Range("P8").Value = 0.9729438112167712
$ws.Range("Q8").Value = 611.1136126944149
$ws.Range("R8").Value = 5500.022514249735
$ws.Range("S8").Value = 0.08603297722130911
$ws.Range("T8").Value = 0.08603297722130909
# Row 9
$ws.Range("I9").Value = 0.08842543241393927
$ws.Range("J9").Value = 0.08842543241393927
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.126567
$ws.Range("N9").Value = 0.379701
$ws.Range("O9").Value = 0.015562836683276
$ws.Range("P9").Value = 0.015562836683276
$ws.Range("Q9").Value = 9.775139365340999
$ws.Range("R9").Value = 87.97625428806899
$ws.Range("S9").Value = 0.001376150563306197
$ws.Range("T9").Value = 0.001376150563306197
# Row 10
$ws.Range("G10").Value = 174.3107043333333
$ws.Range("H10").Value = 522.932113
$ws.Range("I10").Value = 0.199571617988009
$ws.Range("J10").Value = 0.199571617988009
$ws.Range("O10").Value = 0.001498364820294181
$ws.Range("P10").Value = 0.001498364820294181
$ws.Range("Q10").Value = 2.124092139437888
$ws.Range("R10").Value = 19.116829254941
$ws.Range("S10").Value = 0.0002990310915224221
$ws.Range("T10").Value = 0.000299031091522422
# Row 11
$ws.Range("G11").Value = 174.3107043333333
$ws.Range("H11").Value = 522.932113
$ws.Range("I11").Value = 0.199571617988009
$ws.Range("J11").Value = 0.199571617988009
$ws.Range("M11").Value = 0.08128566666666666
$ws.Range("N11").Value = 0.243857
$ws.Range("O11").Value = 0.009994987279658562
$ws.Range("P11").Value = 0.009994987279658561
$ws.Range("Q11").Value = 14.16896180887122
$ws.Range("R11").Value = 127.520656279841
$ws.Range("S11").Value = 0.001994715783171028
$ws.Range("T11").Value = 0.001994715783171028
# Row 12
$ws.Range("G12").Value = 174.3107043333333
$ws.Range("H12").Value = 522.932113
$ws.Range("I12").Value = 0.199571617988009
$ws.Range("J12").Value = 0.199571617988009
$ws.Range("M12").Value = 7.912604999999999
$ws.Range("N12").Value = 23.737815
$ws.Range("O12").Value = 0.9729438112167713
$ws.Range("P12").Value = 0.9729438112167712
$ws.Range("Q12").Value = 1379.251750661455
$ws.Range("R12").Value = 12413.26575595309
$ws.Range("S12").Value = 0.1941719706159511
$ws.Range("T12").Value = 0.194171970615951
# Row 13
$ws.Range("G13").Value = 174.3107043333333
$ws.Range("H13").Value = 522.932113
$ws.Range("I13").Value = 0.199571617988009
$ws.Range("J13").Value = 0.199571617988009
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.126567
$ws.Range("N13").Value = 0.379701
$ws.Range("O13").Value = 0.015562836683276
$ws.Range("P13").Value = 0.015562836683276
$ws.Range("Q13").Value = 22.061982915357
$ws.Range("R13").Value = 198.557846238213
$ws.Range("S13").Value = 0.003105900497364531
$ws.Range("T13").Value = 0.003105900497364531
# Row 14
$ws.Range("G14").Value = 28.53474833333333
$ws.Range("H14").Value = 85.60424499999999
$ws.Range("I14").Value = 0.03266997236655063
$ws.Range("J14").Value = 0.03266997236655063
$ws.Range("O14").Value = 0.001498364820294181
$ws.Range("P14").Value = 0.001498364820294181
$ws.Range("Q14").Value = 0.3477149316072222
$ws.Range("R14").Value = 3.129434384465
$ws.Range("S14").Value = 0.00004895153727402249
$ws.Range("T14").Value = 0.00004895153727402249
# Row 15
$ws.Range("G15").Value = 28.53474833333333
$ws.Range("H15").Value = 85.60424499999999
$ws.Range("I15").Value = 0.03266997236655063
$ws.Range("J15").Value = 0.03266997236655063
$ws.Range("M15").Value = 0.08128566666666666
$ws.Range("N15").Value = 0.243857
$ws.Range("O15").Value = 0.009994987279658562
$ws.Range("P15").Value = 0.009994987279658561
$ws.Range("Q15").Value = 2.319466041440555
$ws.Range("R15").Value = 20.875194372965
$ws.Range("S15").Value = 0.0003265359582304703
$ws.Range("T15").Value = 0.0003265359582304702
# Row 16
$ws.Range("G16").Value = 28.53474833333333
$ws.Range("H16").Value = 85.60424499999999
$ws.Range("I16").Value = 0.03266997236655063
$ws.Range("J16").Value = 0.03266997236655063
$ws.Range("M16").Value = 7.912604999999999
$ws.Range("N16").Value = 23.737815
$ws.Range("O16").Value = 0.9729438112167713
$ws.Range("P16").Value = 0.9729438112167712
$ws.Range("Q16").Value = 225.7841923360749
$ws.Range("R16").Value = 2032.057731024675
$ws.Range("S16").Value = 0.03178604742665837
$ws.Range("T16").Value = 0.03178604742665837
# Row 17
$ws.Range("G17").Value = 28.53474833333333
$ws.Range("H17").Value = 85.60424499999999
$ws.Range("I17").Value = 0.03266997236655063
$ws.Range("J17").Value = 0.03266997236655063
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.126567
$ws.Range("N17").Value = 0.379701
$ws.Range("O17").Value = 0.015562836683276
$ws.Range("P17").Value = 0.015562836683276
$ws.Range("Q17").Value = 3.611557492304999
$ws.Range("R17").Value = 32.50401743074499
$ws.Range("S17").Value = 0.0005084374443877674
$ws.Range("T17").Value = 0.0005084374443877673
# Row 18
$ws.Range("G18").Value = 230.32901
$ws.Range("H18").Value = 690.98703
$ws.Range("I18").Value = 0.263708034289011
$ws.Range("J18").Value = 0.263708034289011
$ws.Range("O18").Value = 0.001498364820294181
$ws.Range("P18").Value = 0.001498364820294181
$ws.Range("Q18").Value = 2.806712539523333
$ws.Range("R18").Value = 25.26041285571
$ws.Range("S18").Value = 0.0003951308414075856
$ws.Range("T18").Value = 0.0003951308414075856
# Row 19
$ws.Range("G19").Value = 230.32901
$ws.Range("H19").Value = 690.98703
$ws.Range("I19").Value = 0.263708034289011
$ws.Range("J19").Value = 0.263708034289011
$ws.Range("M19").Value = 0.08128566666666666
$ws.Range("N19").Value = 0.243857
$ws.Range("O19").Value = 0.009994987279658562
$ws.Range("P19").Value = 0.009994987279658561
$ws.Range("Q19").Value = 18.72244713052333
$ws.Range("R19").Value = 168.50202417471
$ws.Range("S19").Value = 0.002635758448262429
$ws.Range("T19").Value = 0.002635758448262428
# Row 20
$ws.Range("G20").Value = 230.32901
$ws.Range("H20").Value = 690.98703
$ws.Range("I20").Value = 0.263708034289011
$ws.Range("J20").Value = 0.263708034289011
$ws.Range("M20").Value = 7.912604999999999
$ws.Range("N20").Value = 23.737815
$ws.Range("O20").Value = 0.9729438112167713
$ws.Range("P20").Value = 0.9729438112167712
$ws.Range("Q20").Value = 1822.50247617105
$ws.Range("R20").Value = 16402.52228553945
$ws.Range("S20").Value = 0.2565730999296333
$ws.Range("T20").Value = 0.2565730999296333
# Row 21
$ws.Range("G21").Value = 230.32901
$ws.Range("H21").Value = 690.98703
$ws.Range("I21").Value = 0.263708034289011
$ws.Range("J21").Value = 0.263708034289011
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.126567
$ws.Range("N21").Value = 0.379701
$ws.Range("O21").Value = 0.015562836683276
$ws.Range("P21").Value = 0.015562836683276
$ws.Range("Q21").Value = 29.15205180867
$ws.Range("R21").Value = 262.36846627803
$ws.Range("S21").Value = 0.004104045069707625
$ws.Range("T21").Value = 0.004104045069707625
# Row 22
$ws.Range("G22").Value = 151.5554656666667
$ws.Range("H22").Value = 454.666397
$ws.Range("I22").Value = 0.1735187154093718
$ws.Range("J22").Value = 0.1735187154093718
$ws.Range("O22").Value = 0.001498364820294181
$ws.Range("P22").Value = 0.001498364820294181
$ws.Range("Q22").Value = 1.846804386125444
$ws.Range("R22").Value = 16.621239475129
$ws.Range("S22").Value = 0.0002599943388320405
$ws.Range("T22").Value = 0.0002599943388320405
# Row 23
$ws.Range("G23").Value = 151.5554656666667
$ws.Range("H23").Value = 454.666397
$ws.Range("I23").Value = 0.1735187154093718
$ws.Range("J23").Value = 0.1735187154093718
$ws.Range("M23").Value = 0.08128566666666666
$ws.Range("N23").Value = 0.243857
$ws.Range("O23").Value = 0.009994987279658562
$ws.Range("P23").Value = 0.009994987279658561
$ws.Range("Q23").Value = 12.31928706369211
$ws.Range("R23").Value = 110.873583573229
$ws.Range("S23").Value = 0.001734317353299365
$ws.Range("T23").Value = 0.001734317353299365
# Row 24
$ws.Range("G24").Value = 151.5554656666667
$ws.Range("H24").Value = 454.666397
$ws.Range("I24").Value = 0.1735187154093718
$ws.Range("J24").Value = 0.1735187154093718
$ws.Range("M24").Value = 7.912604999999999
$ws.Range("N24").Value = 23.737815
$ws.Range("O24").Value = 0.9729438112167713
$ws.Range("P24").Value = 0.9729438112167712
$ws.Range("Q24").Value = 1199.198535411395
$ws.Range("R24").Value = 10792.78681870255
$ws.Range("S24").Value = 0.1688239602878325
$ws.Range("T24").Value = 0.1688239602878325
# Row 25
$ws.Range("G25").Value = 151.5554656666667
$ws.Range("H25").Value = 454.666397
$ws.Range("I25").Value = 0.1735187154093718
$ws.Range("J25").Value = 0.1735187154093718
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.126567
$ws.Range("N25").Value = 0.379701
$ws.Range("O25").Value = 0.015562836683276
$ws.Range("P25").Value = 0.015562836683276
$ws.Range("Q25").Value = 19.181920623033
$ws.Range("R25").Value = 172.637285607297
$ws.Range("S25").Value = 0.0027004434294079
$ws.Range("T25").Value = 0.0027004434294079
